$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 33-35 (cyclic reshuffle of match data, columns F:V) ---
# Row 33
$ws.Range("F33").Value = "Metz"
$ws.Range("G33").Value = 2
$ws.Range("H33").Value = "Reims"
$ws.Range("I33").Value = 2
$ws.Range("J33").Value = 3.59
$ws.Range("K33").Value = "28/08/2023 05:44"
$ws.Range("L33").Value = 4.82
$ws.Range("M33").Value = "03/09/2023 14:58"
$ws.Range("N33").Value = 3.6
$ws.Range("O33").Value = "28/08/2023 05:44"
$ws.Range("P33").Value = 3.92
$ws.Range("Q33").Value = "03/09/2023 14:57"
$ws.Range("R33").Value = 2.09
$ws.Range("S33").Value = "28/08/2023 05:44"
$ws.Range("T33").Value = 1.78
$ws.Range("U33").Value = "03/09/2023 14:58"
$ws.Range("V33").Value = "https://www.betexplorer.com/football/france/ligue-1/metz-reims/xMcWr6ls/"

# Row 34
$ws.Range("F34").Value = "Le Havre"
$ws.Range("G34").Value = 3
$ws.Range("H34").Value = "Lorient"
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 2.74
$ws.Range("K34").Value = "24/08/2023 07:58"
$ws.Range("L34").Value = 2.68
$ws.Range("M34").Value = "03/09/2023 14:48"
$ws.Range("N34").Value = 3.27
$ws.Range("O34").Value = "24/08/2023 07:58"
$ws.Range("P34").Value = 3.09
$ws.Range("Q34").Value = "03/09/2023 14:59"
$ws.Range("R34").Value = 2.6
$ws.Range("S34").Value = "24/08/2023 07:58"
$ws.Range("T34").Value = 3.08
$ws.Range("U34").Value = "03/09/2023 14:56"
$ws.Range("V34").Value = "https://www.betexplorer.com/football/france/ligue-1/le-havre-lorient/rJv9R4J6/"

# Row 35
$ws.Range("F35").Value = "Lille"
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = "Montpellier"
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 1.51
$ws.Range("K35").Value = "20/08/2023 09:02"
$ws.Range("L35").Value = 1.92
$ws.Range("M35").Value = "03/09/2023 14:58"
$ws.Range("N35").Value = 4.67
$ws.Range("O35").Value = "20/08/2023 09:02"
$ws.Range("P35").Value = 3.92
$ws.Range("Q35").Value = "03/09/2023 14:57"
$ws.Range("R35").Value = 5.42
$ws.Range("S35").Value = "20/08/2023 09:02"
$ws.Range("T35").Value = 4.07
$ws.Range("U35").Value = "03/09/2023 14:58"
$ws.Range("V35").Value = "https://www.betexplorer.com/football/france/ligue-1/lille-montpellier/4YSHPrlJ/"

# --- Update rows 42-44 (cyclic reshuffle of match data, columns F:V) ---
# Row 42
$ws.Range("F42").Value = "Strasbourg"
$ws.Range("G42").Value = 2
$ws.Range("H42").Value = "Montpellier"
$ws.Range("I42").Value = 2
$ws.Range("J42").Value = 2.02
$ws.Range("K42").Value = "28/08/2023 16:01"
$ws.Range("L42").Value = 3.1
$ws.Range("M42").Value = "17/09/2023 14:58"
$ws.Range("N42").Value = 3.61
$ws.Range("O42").Value = "28/08/2023 16:01"
$ws.Range("P42").Value = 3.37
$ws.Range("Q42").Value = "17/09/2023 14:31"
$ws.Range("R42").Value = 3.8
$ws.Range("S42").Value = "28/08/2023 16:01"
$ws.Range("T42").Value = 2.44
$ws.Range("U42").Value = "17/09/2023 14:57"
$ws.Range("V42").Value = "https://www.betexplorer.com/football/france/ligue-1/strasbourg-montpellier/fJq2dPIt/"

# Row 43
$ws.Range("F43").Value = "Clermont"
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = "Nantes"
$ws.Range("I43").Value = 1
$ws.Range("J43").Value = 2.18
$ws.Range("K43").Value = "28/08/2023 16:01"
$ws.Range("L43").Value = 1.88
$ws.Range("M43").Value = "17/09/2023 14:58"
$ws.Range("N43").Value = 3.38
$ws.Range("O43").Value = "28/08/2023 16:01"
$ws.Range("P43").Value = 3.8
$ws.Range("Q43").Value = "17/09/2023 14:58"
$ws.Range("R43").Value = 3.59
$ws.Range("S43").Value = "28/08/2023 16:01"
$ws.Range("T43").Value = 4.34
$ws.Range("U43").Value = "17/09/2023 14:58"
$ws.Range("V43").Value = "https://www.betexplorer.com/football/france/ligue-1/clermont-nantes/Ox0rt4Ya/"

# Row 44
$ws.Range("F44").Value = "Reims"
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = "Brest"
$ws.Range("I44").Value = 2
$ws.Range("J44").Value = 1.71
$ws.Range("K44").Value = "28/08/2023 16:01"
$ws.Range("L44").Value = 2.03
$ws.Range("M44").Value = "17/09/2023 14:50"
$ws.Range("N44").Value = 3.98
$ws.Range("O44").Value = "28/08/2023 16:01"
$ws.Range("P44").Value = 3.6
$ws.Range("Q44").Value = "17/09/2023 14:53"
$ws.Range("R44").Value = 5.01
$ws.Range("S44").Value = "28/08/2023 16:01"
$ws.Range("T44").Value = 3.96
$ws.Range("U44").Value = "17/09/2023 14:53"
$ws.Range("V44").Value = "https://www.betexplorer.com/football/france/ligue-1/reims-brest/pn1vspJg/"

# --- Add new row 47 (copy style from row 46, then set values) ---
$ws.Range("A46").Copy()
$ws.Range("A47").PasteSpecial(-4122)
$ws.Range("E46").Copy()
$ws.Range("E47").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A47").Value = 46
$ws.Range("B47").Value = "france"
$ws.Range("C47").Value = "ligue-1"
$ws.Range("D47").Value = "2023-2024"
$ws.Range("E47").Value = 45191.875
$ws.Range("F47").Value = "Monaco"
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = "Nice"
$ws.Range("I47").Value = 1
$ws.Range("J47").Value = 1.85
$ws.Range("K47").Value = "05/09/2023 12:01"
$ws.Range("L47").Value = 2.46
$ws.Range("M47").Value = "22/09/2023 20:58"
$ws.Range("N47").Value = 3.93
$ws.Range("O47").Value = "05/09/2023 12:01"
$ws.Range("P47").Value = 3.63
$ws.Range("Q47").Value = "22/09/2023 20:40"
$ws.Range("R47").Value = 4.15
$ws.Range("S47").Value = "05/09/2023 12:01"
$ws.Range("T47").Value = 2.95
$ws.Range("U47").Value = "22/09/2023 20:58"
$ws.Range("V47").Value = "https://www.betexplorer.com/football/france/ligue-1/monaco-nice/ImqztTyg/"
